$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format first so numeric-looking values
# (e.g. "399.19") are written as strings, not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.353.45"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").Value = "3.346.33"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "399.19"
$ws.Range("E5").Value = "  -3.60%  "

$ws.Range("D6").Value = "125.33"
$ws.Range("E6").Value = "  +7.91%  "

$ws.Range("E7").Value = "  +2.00%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").Value = "0.656"
$ws.Range("E9").Value = "  +4.19%  "

$ws.Range("E10").Value = "  +0.30%  "

$ws.Range("D11").Value = "40.60"
$ws.Range("E11").Value = "  +1.74%  "

$ws.Range("E12").Value = "  -1.01%  "

$ws.Range("D13").Value = "3.866.83"
$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("D14").Value = "8.22"
$ws.Range("E14").Value = "  -1.24%  "

$ws.Range("D15").Value = "19.19"
$ws.Range("E15").Value = "  -0.32%  "

$ws.Range("D16").Value = "3.327.97"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").Value = "61.240.16"
$ws.Range("E17").Value = "  -0.84%  "

$ws.Range("D18").Value = "11.13"
$ws.Range("E18").Value = "  +2.54%  "

$ws.Range("D19").Value = "1.00"
$ws.Range("E19").Value = "  -0.71%  "

$ws.Range("D20").Value = "0.0000125"
$ws.Range("E20").Value = "  +5.83%  "

$ws.Range("D21").Value = "3.19"
$ws.Range("E21").Value = "  -4.22%  "

$ws.Range("D22").Value = "79.59"
$ws.Range("E22").Value = "  +6.56%  "

$ws.Range("D23").Value = "12.66"
$ws.Range("E23").Value = "  +0.91%  "

$ws.Range("D24").Value = "297.90"
$ws.Range("E24").Value = "  +0.98%  "

$ws.Range("D25").Value = "3.08"
$ws.Range("E25").Value = "  -1.66%  "

$ws.Range("D26").Value = "4.75"
$ws.Range("E26").Value = "  +11.47%  "

$ws.Range("D27").Value = "28.87"
$ws.Range("E27").Value = "  -1.85%  "

$ws.Range("D28").Value = "8.13"
$ws.Range("E28").Value = "  +6.65%  "

$ws.Range("D29").Value = "7.42"
$ws.Range("E29").Value = "  -6.38%  "

$ws.Range("E30").Value = "  -1.82%  "

$ws.Range("D31").Value = "0.113"
$ws.Range("E31").Value = "  -0.55%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").Value = "11.28"
$ws.Range("E33").Value = "  -1.50%  "

$ws.Range("D34").Value = "2.50"
$ws.Range("E34").Value = "  -3.37%  "

$ws.Range("D35").Value = "40.90"
$ws.Range("E35").Value = "  -4.43%  "

$ws.Range("D36").Value = "0.0476"
$ws.Range("E36").Value = "  -3.13%  "

$ws.Range("D37").Value = "51.94"
$ws.Range("E37").Value = "  -0.62%  "

$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").Value = "3.35"
$ws.Range("E39").Value = "  -2.72%  "

$ws.Range("D40").Value = "2.89"
$ws.Range("E40").Value = "  -7.98%  "

$ws.Range("D41").Value = "136.88"
$ws.Range("E41").Value = "  +2.54%  "

$ws.Range("E42").Value = "  +2.18%  "

$ws.Range("D43").Value = "0.122"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.279"
$ws.Range("E44").Value = "  -2.39%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "3.85"
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("D46").Value = "16.48"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "2.23"
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("D48").Value = "20.93"
$ws.Range("E48").Value = "  -1.14%  "

$ws.Range("D49").Value = "3.671.93"
$ws.Range("E49").Value = "  -0.17%  "

$ws.Range("D50").Value = "2.088.99"
$ws.Range("E50").Value = "  -3.80%  "

$ws.Range("D51").Value = "2.29"
$ws.Range("E51").Value = "  -5.04%  "

# Clear the temporary text formatting so no stray style survives on column D,
# restoring cells to their original (unstyled) appearance.
$ws.Range("D2:D51").ClearFormats()

Write-Host "Updated cryptos list"
